$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B17's number format to match the datetime format (style s=2) instead of date-only (s=3)
$ws.Range("B17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 18 with data for June 17th (date serial 43999)
# First, copy formatting from row 17 so borders/fonts/number formats line up,
# then overwrite the values that differ.
$ws.Range("A17:G17").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 43999
$ws.Range("B18").NumberFormat = "YYYY-MM-DD"
$ws.Range("C18").Value = 159793
$ws.Range("D18").Value = 222801
$ws.Range("E18").Value = 59076
$ws.Range("F18").Value = 19080
$ws.Range("G18").Value = 32.06
